# Range bar chart added
# Adds the MAX/MIN "range" helper cells in column N that feed the new
# range bar chart, then leaves the selection where Excel would after
# typing them in (bottom-right new cell, N24 view scrolled to A3).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# PC Gamers block (rows 5:6)
$ws.Range("N5").Formula = "=MAX(C5:G6)"
$ws.Range("N6").Formula = "=MIN(C5:G6)"

# Console Gamers block (rows 11:14)
$ws.Range("N13").Formula = "=MAX(C11:G14)"
$ws.Range("N14").Formula = "=MIN(C11:G14)"

# Non-Gamers block (rows 19:23)
$ws.Range("N22").Formula = "=MAX(C19:G23)"
$ws.Range("N23").Formula = "=MIN(C19:G23)"

# Match the saved view/selection state from the diff.
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("N24").Select()
